$d = $word.ActiveDocument

# 1) Merge the "cyclonomade" spell-check-wrapped text into the surrounding
#    run text (removing the proofErr markers / run split) by replacing the
#    whole sentence with a single run's worth of text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Je me présente, Max le Bourlingueur, cyclonomade professionnel. Depuis bientôt cinq ans, je parcours les routes de France en totale autonomie avec mon vélo et ma roulotte, accompagné de mes trois fidèles chiens.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Je me présente, Max le Bourlingueur, cyclonomade professionnel. Depuis bientôt cinq ans, je parcours les routes de France en totale autonomie avec mon vélo et ma roulotte, accompagné de mes trois fidèles chiens.",
    2
)

# 2) Add the Siret number as a new, separate run right after "Siret : "
#    (the paragraph keeps its original identity/attributes; we only split
#    the run so the inserted number lives in its own <w:r>).
$siretRange = $d.Content
$siretRange.Find.ClearFormatting()
$siretRange.Find.Execute("Siret : ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$siretXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4425ADF2" w14:textId="7B7D1B5B" w:rsidR="005A22F3" w:rsidRDefault="005A22F3" w:rsidP="005A22F3"><w:r><w:t xml:space="preserve">Siret : </w:t></w:r><w:r><w:t>99927450900010</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$siretRange.InsertXML($siretXml)
